$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap header labels between B1 and C1
$headerB = $ws.Range("B1").Value()
$headerC = $ws.Range("C1").Value()
$ws.Range("B1").Value = $headerC
$ws.Range("C1").Value = $headerB

# Swap data values between column B and column C for rows 2-13
for ($r = 2; $r -le 13; $r++) {
    $cellB = $ws.Cells.Item($r, 2)
    $cellC = $ws.Cells.Item($r, 3)
    $valB = $cellB.Value()
    $valC = $cellC.Value()
    $cellB.Value = $valC
    $cellC.Value = $valB
}
